$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 245, pushing the existing
# rows 245:260 down to 247:262 (formatting carries down with them).
$ws.Rows("245:246").Insert()

# Populate the two new rows (245 and 246) with the new records.
# Columns A,B,C,E,F,G,H,I,R are constant across this data block.

# Row 245
$ws.Cells.Item(245, 1).Value = 3
$ws.Cells.Item(245, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(245, 3).Value = "Coquimbo"
$ws.Cells.Item(245, 4).Value = 44516
$ws.Cells.Item(245, 5).Value = 5
$ws.Cells.Item(245, 6).Value = 100112032
$ws.Cells.Item(245, 7).Value = "Zapallo italiano"
$ws.Cells.Item(245, 8).Value = "Sin especificar"
$ws.Cells.Item(245, 9).Value = "Primera"
$ws.Cells.Item(245, 10).Value = 130
$ws.Cells.Item(245, 11).Value = 4000
$ws.Cells.Item(245, 12).Value = 4500
$ws.Cells.Item(245, 13).Value = 4269
$ws.Cells.Item(245, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(245, 15).Value = "Limache"
$ws.Cells.Item(245, 16).Value = 119
$ws.Cells.Item(245, 17).Value = 36
$ws.Cells.Item(245, 18).Value = "Hortaliza"

# Row 246
$ws.Cells.Item(246, 1).Value = 3
$ws.Cells.Item(246, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(246, 3).Value = "Coquimbo"
$ws.Cells.Item(246, 4).Value = 44516
$ws.Cells.Item(246, 5).Value = 5
$ws.Cells.Item(246, 6).Value = 100112032
$ws.Cells.Item(246, 7).Value = "Zapallo italiano"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 125
$ws.Cells.Item(246, 11).Value = 8000
$ws.Cells.Item(246, 12).Value = 8500
$ws.Cells.Item(246, 13).Value = 8260
$ws.Cells.Item(246, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(246, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(246, 16).Value = 118
$ws.Cells.Item(246, 17).Value = 70
$ws.Cells.Item(246, 18).Value = "Hortaliza"
